$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated line-flow simulation results (380 kV case) for rows 2-25, columns B-O
# Each entry is (ColumnLetter, RowNumber, NewValue)
$updates = New-Object System.Collections.ArrayList
[void]$updates.Add(@("B", 2, 0.5727304149500867))
[void]$updates.Add(@("C", 2, 0.07761289738254362))
[void]$updates.Add(@("D", 2, 0.07787587772119764))
[void]$updates.Add(@("E", 2, 0.08989345109232616))
[void]$updates.Add(@("G", 2, 0.0024881659178312))
[void]$updates.Add(@("I", 2, 1.095153994367458))
[void]$updates.Add(@("K", 2, 0.3696633564866261))
[void]$updates.Add(@("L", 2, 0.214151979812641))
[void]$updates.Add(@("O", 2, 4.472142840385061))
[void]$updates.Add(@("B", 3, 0.5327881316985952))
[void]$updates.Add(@("C", 3, 0.07526919031080581))
[void]$updates.Add(@("D", 3, 0.07073728468191121))
[void]$updates.Add(@("E", 3, 0.08938513881929566))
[void]$updates.Add(@("G", 3, 0.002491019998685183))
[void]$updates.Add(@("I", 3, 1.100956501799899))
[void]$updates.Add(@("K", 3, 0.3321063022459043))
[void]$updates.Add(@("L", 3, 0.2069025281903834))
[void]$updates.Add(@("O", 3, 4.485218395355673))
[void]$updates.Add(@("B", 4, 0.5084538199127167))
[void]$updates.Add(@("C", 4, 0.07381036626509285))
[void]$updates.Add(@("D", 4, 0.06638880891722465))
[void]$updates.Add(@("E", 4, 0.08911695983380064))
[void]$updates.Add(@("G", 4, 0.002492866455910748))
[void]$updates.Add(@("I", 4, 1.105031197411417))
[void]$updates.Add(@("K", 4, 0.3091043750108327))
[void]$updates.Add(@("L", 4, 0.2025679952275539))
[void]$updates.Add(@("O", 4, 4.495504722827405))
[void]$updates.Add(@("B", 5, 0.4985858389825069))
[void]$updates.Add(@("C", 5, 0.07321092045347655))
[void]$updates.Add(@("D", 5, 0.06462548077870167))
[void]$updates.Add(@("E", 5, 0.08901873485859824))
[void]$updates.Add(@("G", 5, 0.002493642620709251))
[void]$updates.Add(@("I", 5, 1.10682039297318))
[void]$updates.Add(@("K", 5, 0.2997459926357351))
[void]$updates.Add(@("L", 5, 0.2008310268254405))
[void]$updates.Add(@("O", 5, 4.500264162182702))
[void]$updates.Add(@("B", 6, 0.4969502112424493))
[void]$updates.Add(@("C", 6, 0.07311108347885664))
[void]$updates.Add(@("D", 6, 0.06433320756431726))
[void]$updates.Add(@("E", 6, 0.08900309310572752))
[void]$updates.Add(@("G", 6, 0.002493772936698238))
[void]$updates.Add(@("I", 6, 1.10712526199503))
[void]$updates.Add(@("K", 6, 0.2981929647971242))
[void]$updates.Add(@("L", 6, 0.2005443811665799))
[void]$updates.Add(@("O", 6, 4.501088747328083))
[void]$updates.Add(@("B", 7, 0.5083205399478743))
[void]$updates.Add(@("C", 7, 0.07380230201050608))
[void]$updates.Add(@("D", 7, 0.06636499278894803))
[void]$updates.Add(@("E", 7, 0.08911559033606764))
[void]$updates.Add(@("G", 7, 0.002492876827318294))
[void]$updates.Add(@("I", 7, 1.105054805939496))
[void]$updates.Add(@("K", 7, 0.3089781029089238))
[void]$updates.Add(@("L", 7, 0.2025444507974328))
[void]$updates.Add(@("O", 7, 4.495566611910647))
[void]$updates.Add(@("B", 8, 0.558919122507632))
[void]$updates.Add(@("C", 8, 0.07680889770454513))
[void]$updates.Add(@("D", 8, 0.07540729610471431))
[void]$updates.Add(@("E", 8, 0.08970907626408575))
[void]$updates.Add(@("G", 8, 0.002489130529754434))
[void]$updates.Add(@("I", 8, 1.097048442995174))
[void]$updates.Add(@("K", 8, 0.3567018450827391))
[void]$updates.Add(@("L", 8, 0.2116281857538667))
[void]$updates.Add(@("O", 8, 4.476182608521242))
[void]$updates.Add(@("B", 9, 0.6596341010568665))
[void]$updates.Add(@("C", 9, 0.08254768830322945))
[void]$updates.Add(@("D", 9, 0.0934151888224477))
[void]$updates.Add(@("E", 9, 0.09122097426681108))
[void]$updates.Add(@("G", 9, 0.002482526969898838))
[void]$updates.Add(@("I", 9, 1.08541074324107))
[void]$updates.Add(@("K", 9, 0.4507356643472349))
[void]$updates.Add(@("L", 9, 0.2303658104416826))
[void]$updates.Add(@("O", 9, 4.456093704399052))
[void]$updates.Add(@("B", 10, 0.7345203339071986))
[void]$updates.Add(@("C", 10, 0.08666828313526764))
[void]$updates.Add(@("D", 10, 0.1068166363773173))
[void]$updates.Add(@("E", 10, 0.09254358144967156))
[void]$updates.Add(@("G", 10, 0.002478123684613767))
[void]$updates.Add(@("I", 10, 1.07933918677022))
[void]$updates.Add(@("K", 10, 0.5200828564707649))
[void]$updates.Add(@("L", 10, 0.2446960477828242))
[void]$updates.Add(@("O", 10, 4.452277173103852))
[void]$updates.Add(@("B", 11, 0.7687781501575159))
[void]$updates.Add(@("C", 11, 0.0885221148368629))
[void]$updates.Add(@("D", 11, 0.1129511222191724))
[void]$updates.Add(@("E", 11, 0.09319119891983618))
[void]$updates.Add(@("G", 11, 0.002476216902647413))
[void]$updates.Add(@("I", 11, 1.077115772420065))
[void]$updates.Add(@("K", 11, 0.5516851078120908))
[void]$updates.Add(@("L", 11, 0.2513377848951137))
[void]$updates.Add(@("O", 11, 4.45292099780562))
[void]$updates.Add(@("B", 12, 0.7817778045026671))
[void]$updates.Add(@("C", 12, 0.08922113882944416))
[void]$updates.Add(@("D", 12, 0.1152795934364832))
[void]$updates.Add(@("E", 12, 0.09344303284816036))
[void]$updates.Add(@("G", 12, 0.002475508627789458))
[void]$updates.Add(@("I", 12, 1.076351298516094))
[void]$updates.Add(@("K", 12, 0.5636597521410636))
[void]$updates.Add(@("L", 12, 0.2538704744977451))
[void]$updates.Add(@("O", 12, 4.453507265683442))
[void]$updates.Add(@("B", 13, 0.7789769067306054))
[void]$updates.Add(@("C", 13, 0.08907072432511143))
[void]$updates.Add(@("D", 13, 0.11477787208797))
[void]$updates.Add(@("E", 13, 0.09338850279245747))
[void]$updates.Add(@("G", 13, 0.00247566055557185))
[void]$updates.Add(@("I", 13, 1.076512494858626))
[void]$updates.Add(@("K", 13, 0.5610804696561615))
[void]$updates.Add(@("L", 13, 0.2533242322596294))
[void]$updates.Add(@("O", 13, 4.453365767556562))
[void]$updates.Add(@("B", 14, 0.7698471026100719))
[void]$updates.Add(@("C", 14, 0.08857968379429337))
[void]$updates.Add(@("D", 14, 0.1131425772544787))
[void]$updates.Add(@("E", 14, 0.09321178534632324))
[void]$updates.Add(@("G", 14, 0.002476158356603485))
[void]$updates.Add(@("I", 14, 1.077051325553342))
[void]$updates.Add(@("K", 14, 0.5526701202595063))
[void]$updates.Add(@("L", 14, 0.2515457982278946))
[void]$updates.Add(@("O", 14, 4.452962365602929))
[void]$updates.Add(@("B", 15, 0.7642583274908645))
[void]$updates.Add(@("C", 15, 0.08827851859174984))
[void]$updates.Add(@("D", 15, 0.1121416247675313))
[void]$updates.Add(@("E", 15, 0.09310439930889558))
[void]$updates.Add(@("G", 15, 0.002476465067012852))
[void]$updates.Add(@("I", 15, 1.077391466902483))
[void]$updates.Add(@("K", 15, 0.5475195085508915))
[void]$updates.Add(@("L", 15, 0.250458747425057))
[void]$updates.Add(@("O", 15, 4.452759876086134))
[void]$updates.Add(@("B", 16, 0.7322853108473169))
[void]$updates.Add(@("C", 16, 0.08654671430106475))
[void]$updates.Add(@("D", 16, 0.1064164996273433))
[void]$updates.Add(@("E", 16, 0.09250218175887426))
[void]$updates.Add(@("G", 16, 0.00247825022938034))
[void]$updates.Add(@("I", 16, 1.07949533285452))
[void]$updates.Add(@("K", 16, 0.5180186584584874))
[void]$updates.Add(@("L", 16, 0.2442644616873366))
[void]$updates.Add(@("O", 16, 4.452283008388605))
[void]$updates.Add(@("B", 17, 0.7127195684890353))
[void]$updates.Add(@("C", 17, 0.08547901100675404))
[void]$updates.Add(@("D", 17, 0.1029140763440921))
[void]$updates.Add(@("E", 17, 0.09214450200122926))
[void]$updates.Add(@("G", 17, 0.002479369985887842))
[void]$updates.Add(@("I", 17, 1.080923954565392))
[void]$updates.Add(@("K", 17, 0.4999348117324871))
[void]$updates.Add(@("L", 17, 0.2404958902951506))
[void]$updates.Add(@("O", 17, 4.452600197878212))
[void]$updates.Add(@("B", 18, 0.701483949902439))
[void]$updates.Add(@("C", 18, 0.08486295361859675))
[void]$updates.Add(@("D", 18, 0.1009031615599554))
[void]$updates.Add(@("E", 18, 0.09194310053677412))
[void]$updates.Add(@("G", 18, 0.002480023107493833))
[void]$updates.Add(@("I", 18, 1.081796350099154))
[void]$updates.Add(@("K", 18, 0.4895387480480338))
[void]$updates.Add(@("L", 18, 0.238339873396626))
[void]$updates.Add(@("O", 18, 4.453006638673969))
[void]$updates.Add(@("B", 19, 0.6976828859768602))
[void]$updates.Add(@("C", 19, 0.08465403390407289))
[void]$updates.Add(@("D", 19, 0.1002229164602682))
[void]$updates.Add(@("E", 19, 0.09187565294132227))
[void]$updates.Add(@("G", 19, 0.00248024580261802))
[void]$updates.Add(@("I", 19, 1.082100433111854))
[void]$updates.Add(@("K", 19, 0.4860197420377403))
[void]$updates.Add(@("L", 19, 0.2376118712531792))
[void]$updates.Add(@("O", 19, 4.453182718170808))
[void]$updates.Add(@("B", 20, 0.7148005084369231))
[void]$updates.Add(@("C", 20, 0.08559287099966184))
[void]$updates.Add(@("D", 20, 0.1032865441273145))
[void]$updates.Add(@("E", 20, 0.09218212994475294))
[void]$updates.Add(@("G", 20, 0.002479249847847379))
[void]$updates.Add(@("I", 20, 1.080766628820271))
[void]$updates.Add(@("K", 20, 0.5018593243942178))
[void]$updates.Add(@("L", 20, 0.2408958645564212))
[void]$updates.Add(@("O", 20, 4.452543247517156))
[void]$updates.Add(@("B", 21, 0.7725280204149954))
[void]$updates.Add(@("C", 21, 0.08872399532883435))
[void]$updates.Add(@("D", 21, 0.1136227544458421))
[void]$updates.Add(@("E", 21, 0.09326351265699273))
[void]$updates.Add(@("G", 21, 0.002476011766779251))
[void]$updates.Add(@("I", 21, 1.076890954684416))
[void]$updates.Add(@("K", 21, 0.5551402424167975))
[void]$updates.Add(@("L", 21, 0.2520676899893601))
[void]$updates.Add(@("O", 21, 4.453071558316736))
[void]$updates.Add(@("B", 22, 0.8104131689062228))
[void]$updates.Add(@("C", 22, 0.09075297983991959))
[void]$updates.Add(@("D", 22, 0.1204099730264545))
[void]$updates.Add(@("E", 22, 0.09400869472696627))
[void]$updates.Add(@("G", 22, 0.002473975796212133))
[void]$updates.Add(@("I", 22, 1.074809637041383))
[void]$updates.Add(@("K", 22, 0.5900062347294579))
[void]$updates.Add(@("L", 22, 0.2594717297842379))
[void]$updates.Add(@("O", 22, 4.455413095566797))
[void]$updates.Add(@("B", 23, 0.7901789957709298))
[void]$updates.Add(@("C", 23, 0.08967166763281398))
[void]$updates.Add(@("D", 23, 0.1167845887102033))
[void]$updates.Add(@("E", 23, 0.093607464539879))
[void]$updates.Add(@("G", 23, 0.002475055105703543))
[void]$updates.Add(@("I", 23, 1.075879135482971))
[void]$updates.Add(@("K", 23, 0.5713937444136548))
[void]$updates.Add(@("L", 23, 0.2555106846090638))
[void]$updates.Add(@("O", 23, 4.453980640763945))
[void]$updates.Add(@("B", 24, 0.7138596749886688))
[void]$updates.Add(@("C", 24, 0.08554140181287551))
[void]$updates.Add(@("D", 24, 0.1031181431109331))
[void]$updates.Add(@("E", 24, 0.09216510516427689))
[void]$updates.Add(@("G", 24, 0.00247930413313981))
[void]$updates.Add(@("I", 24, 1.080837596734547))
[void]$updates.Add(@("K", 24, 0.5009892503770743))
[void]$updates.Add(@("L", 24, 0.2407150032190231))
[void]$updates.Add(@("O", 24, 4.45256829676876))
[void]$updates.Add(@("B", 25, 0.632230268845376))
[void]$updates.Add(@("C", 25, 0.08101200613673143))
[void]$updates.Add(@("D", 25, 0.08851374018803426))
[void]$updates.Add(@("E", 25, 0.09077474121965423))
[void]$updates.Add(@("G", 25, 0.002484234345708825))
[void]$updates.Add(@("I", 25, 1.088123908132964))
[void]$updates.Add(@("K", 25, 0.4507356643472349))
[void]$updates.Add(@("L", 25, 0.2251978035643845))
[void]$updates.Add(@("O", 25, 4.459607802651078))

foreach ($u in $updates) {
    $col = $u[0]
    $row = $u[1]
    $val = $u[2]
    $ws.Range("$col$row").Value = $val
}
